# Auto-generated edit script for Daily Report update: 2026-02-23
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Daily_Data - append 24 new rows (218-241) for date 46073 ----
$ws1 = $wb.Worksheets.Item("Daily_Data")

$ws1.Cells.Item(218, 1).Value = 46073
$ws1.Cells.Item(218, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(218, 2).Value = "ASAHI DEPOSITORY LLC Registered"
$ws1.Cells.Item(218, 3).Value = 23301775.992
$ws1.Cells.Item(218, 4).Value = 0
$ws1.Cells.Item(218, 5).Value = 0
$ws1.Cells.Item(218, 6).Value = 0
$ws1.Cells.Item(218, 7).Value = 0
$ws1.Cells.Item(218, 8).Value = 23301775.992

$ws1.Cells.Item(219, 1).Value = 46073
$ws1.Cells.Item(219, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(219, 2).Value = "ASAHI DEPOSITORY LLC Eligible"
$ws1.Cells.Item(219, 3).Value = 2748893.808
$ws1.Cells.Item(219, 4).Value = 0
$ws1.Cells.Item(219, 5).Value = 0
$ws1.Cells.Item(219, 6).Value = 0
$ws1.Cells.Item(219, 7).Value = 0
$ws1.Cells.Item(219, 8).Value = 2748893.808

$ws1.Cells.Item(220, 1).Value = 46073
$ws1.Cells.Item(220, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(220, 2).Value = "BRINK'S, INC. Registered"
$ws1.Cells.Item(220, 3).Value = 15782712.636
$ws1.Cells.Item(220, 4).Value = 0
$ws1.Cells.Item(220, 5).Value = 0
$ws1.Cells.Item(220, 6).Value = 0
$ws1.Cells.Item(220, 7).Value = 0
$ws1.Cells.Item(220, 8).Value = 15782712.636

$ws1.Cells.Item(221, 1).Value = 46073
$ws1.Cells.Item(221, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(221, 2).Value = "BRINK'S, INC. Eligible"
$ws1.Cells.Item(221, 3).Value = 39336942.517
$ws1.Cells.Item(221, 4).Value = 0
$ws1.Cells.Item(221, 5).Value = 0
$ws1.Cells.Item(221, 6).Value = 0
$ws1.Cells.Item(221, 7).Value = 0
$ws1.Cells.Item(221, 8).Value = 39336942.517

$ws1.Cells.Item(222, 1).Value = 46073
$ws1.Cells.Item(222, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(222, 2).Value = "CNT DEPOSITORY, INC. Registered"
$ws1.Cells.Item(222, 3).Value = 12174851.569
$ws1.Cells.Item(222, 4).Value = 0
$ws1.Cells.Item(222, 5).Value = 0
$ws1.Cells.Item(222, 6).Value = 0
$ws1.Cells.Item(222, 7).Value = 0
$ws1.Cells.Item(222, 8).Value = 12174851.569

$ws1.Cells.Item(223, 1).Value = 46073
$ws1.Cells.Item(223, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(223, 2).Value = "CNT DEPOSITORY, INC. Eligible"
$ws1.Cells.Item(223, 3).Value = 14018899.428
$ws1.Cells.Item(223, 4).Value = 0
$ws1.Cells.Item(223, 5).Value = 162211.605
$ws1.Cells.Item(223, 6).Value = -162211.605
$ws1.Cells.Item(223, 7).Value = 0
$ws1.Cells.Item(223, 8).Value = 13856687.823

$ws1.Cells.Item(224, 1).Value = 46073
$ws1.Cells.Item(224, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(224, 2).Value = "DELAWARE DEPOSITORY Registered"
$ws1.Cells.Item(224, 3).Value = 1532776.423
$ws1.Cells.Item(224, 4).Value = 0
$ws1.Cells.Item(224, 5).Value = 0
$ws1.Cells.Item(224, 6).Value = 0
$ws1.Cells.Item(224, 7).Value = 0
$ws1.Cells.Item(224, 8).Value = 1532776.423

$ws1.Cells.Item(225, 1).Value = 46073
$ws1.Cells.Item(225, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(225, 2).Value = "DELAWARE DEPOSITORY Eligible"
$ws1.Cells.Item(225, 3).Value = 16272628.85
$ws1.Cells.Item(225, 4).Value = 0
$ws1.Cells.Item(225, 5).Value = 500751.905
$ws1.Cells.Item(225, 6).Value = -500751.905
$ws1.Cells.Item(225, 7).Value = 0
$ws1.Cells.Item(225, 8).Value = 15771876.945

$ws1.Cells.Item(226, 1).Value = 46073
$ws1.Cells.Item(226, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(226, 2).Value = "HSBC BANK, USA Registered"
$ws1.Cells.Item(226, 3).Value = 3412157.57
$ws1.Cells.Item(226, 4).Value = 0
$ws1.Cells.Item(226, 5).Value = 0
$ws1.Cells.Item(226, 6).Value = 0
$ws1.Cells.Item(226, 7).Value = 0
$ws1.Cells.Item(226, 8).Value = 3412157.57

$ws1.Cells.Item(227, 1).Value = 46073
$ws1.Cells.Item(227, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(227, 2).Value = "HSBC BANK, USA Eligible"
$ws1.Cells.Item(227, 3).Value = 19638409.353
$ws1.Cells.Item(227, 4).Value = 0
$ws1.Cells.Item(227, 5).Value = 627200
$ws1.Cells.Item(227, 6).Value = -627200
$ws1.Cells.Item(227, 7).Value = 0
$ws1.Cells.Item(227, 8).Value = 19011209.353

$ws1.Cells.Item(228, 1).Value = 46073
$ws1.Cells.Item(228, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(228, 2).Value = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered"
$ws1.Cells.Item(228, 3).Value = 273789.87
$ws1.Cells.Item(228, 4).Value = 0
$ws1.Cells.Item(228, 5).Value = 0
$ws1.Cells.Item(228, 6).Value = 0
$ws1.Cells.Item(228, 7).Value = 0
$ws1.Cells.Item(228, 8).Value = 273789.87

$ws1.Cells.Item(229, 1).Value = 46073
$ws1.Cells.Item(229, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(229, 2).Value = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible"
$ws1.Cells.Item(229, 3).Value = 3295246.644
$ws1.Cells.Item(229, 4).Value = 0
$ws1.Cells.Item(229, 5).Value = 0
$ws1.Cells.Item(229, 6).Value = 0
$ws1.Cells.Item(229, 7).Value = 0
$ws1.Cells.Item(229, 8).Value = 3295246.644

$ws1.Cells.Item(230, 1).Value = 46073
$ws1.Cells.Item(230, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(230, 2).Value = "JP MORGAN CHASE BANK NA Registered"
$ws1.Cells.Item(230, 3).Value = 12000343.77
$ws1.Cells.Item(230, 4).Value = 0
$ws1.Cells.Item(230, 5).Value = 0
$ws1.Cells.Item(230, 6).Value = 0
$ws1.Cells.Item(230, 7).Value = 0
$ws1.Cells.Item(230, 8).Value = 12000343.77

$ws1.Cells.Item(231, 1).Value = 46073
$ws1.Cells.Item(231, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(231, 2).Value = "JP MORGAN CHASE BANK NA Eligible"
$ws1.Cells.Item(231, 3).Value = 143854408.433
$ws1.Cells.Item(231, 4).Value = 0
$ws1.Cells.Item(231, 5).Value = 963743.1
$ws1.Cells.Item(231, 6).Value = -963743.1
$ws1.Cells.Item(231, 7).Value = 0
$ws1.Cells.Item(231, 8).Value = 142890665.333

$ws1.Cells.Item(232, 1).Value = 46073
$ws1.Cells.Item(232, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(232, 2).Value = "LOOMIS INTERNATIONAL (US) LLC Registered"
$ws1.Cells.Item(232, 3).Value = 6311885.937
$ws1.Cells.Item(232, 4).Value = 0
$ws1.Cells.Item(232, 5).Value = 0
$ws1.Cells.Item(232, 6).Value = 0
$ws1.Cells.Item(232, 7).Value = 0
$ws1.Cells.Item(232, 8).Value = 6311885.937

$ws1.Cells.Item(233, 1).Value = 46073
$ws1.Cells.Item(233, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(233, 2).Value = "LOOMIS INTERNATIONAL (US) LLC Eligible"
$ws1.Cells.Item(233, 3).Value = 24033585.186
$ws1.Cells.Item(233, 4).Value = 0
$ws1.Cells.Item(233, 5).Value = 0
$ws1.Cells.Item(233, 6).Value = 0
$ws1.Cells.Item(233, 7).Value = 0
$ws1.Cells.Item(233, 8).Value = 24033585.186

$ws1.Cells.Item(234, 1).Value = 46073
$ws1.Cells.Item(234, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(234, 2).Value = "MALCA-AMIT ARMORED, INC. Registered"
$ws1.Cells.Item(234, 3).Value = 0
$ws1.Cells.Item(234, 4).Value = 0
$ws1.Cells.Item(234, 5).Value = 0
$ws1.Cells.Item(234, 6).Value = 0
$ws1.Cells.Item(234, 7).Value = 0
$ws1.Cells.Item(234, 8).Value = 0

$ws1.Cells.Item(235, 1).Value = 46073
$ws1.Cells.Item(235, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(235, 2).Value = "MALCA-AMIT ARMORED, INC. Eligible"
$ws1.Cells.Item(235, 3).Value = 0
$ws1.Cells.Item(235, 4).Value = 0
$ws1.Cells.Item(235, 5).Value = 0
$ws1.Cells.Item(235, 6).Value = 0
$ws1.Cells.Item(235, 7).Value = 0
$ws1.Cells.Item(235, 8).Value = 0

$ws1.Cells.Item(236, 1).Value = 46073
$ws1.Cells.Item(236, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(236, 2).Value = "MALCA-AMIT USA, LLC Registered"
$ws1.Cells.Item(236, 3).Value = 949634.064
$ws1.Cells.Item(236, 4).Value = 0
$ws1.Cells.Item(236, 5).Value = 0
$ws1.Cells.Item(236, 6).Value = 0
$ws1.Cells.Item(236, 7).Value = 0
$ws1.Cells.Item(236, 8).Value = 949634.064

$ws1.Cells.Item(237, 1).Value = 46073
$ws1.Cells.Item(237, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(237, 2).Value = "MALCA-AMIT USA, LLC Eligible"
$ws1.Cells.Item(237, 3).Value = 1073898.377
$ws1.Cells.Item(237, 4).Value = 0
$ws1.Cells.Item(237, 5).Value = 0
$ws1.Cells.Item(237, 6).Value = 0
$ws1.Cells.Item(237, 7).Value = 0
$ws1.Cells.Item(237, 8).Value = 1073898.377

$ws1.Cells.Item(238, 1).Value = 46073
$ws1.Cells.Item(238, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(238, 2).Value = "MANFRA, TORDELLA & BROOKES, LLC Registered"
$ws1.Cells.Item(238, 3).Value = 6219630.033
$ws1.Cells.Item(238, 4).Value = 0
$ws1.Cells.Item(238, 5).Value = 0
$ws1.Cells.Item(238, 6).Value = 0
$ws1.Cells.Item(238, 7).Value = 0
$ws1.Cells.Item(238, 8).Value = 6219630.033

$ws1.Cells.Item(239, 1).Value = 46073
$ws1.Cells.Item(239, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(239, 2).Value = "MANFRA, TORDELLA & BROOKES, LLC Eligible"
$ws1.Cells.Item(239, 3).Value = 12256015.907
$ws1.Cells.Item(239, 4).Value = 0
$ws1.Cells.Item(239, 5).Value = 0
$ws1.Cells.Item(239, 6).Value = 0
$ws1.Cells.Item(239, 7).Value = 0
$ws1.Cells.Item(239, 8).Value = 12256015.907

$ws1.Cells.Item(240, 1).Value = 46073
$ws1.Cells.Item(240, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(240, 2).Value = "STONEX PRECIOUS METALS LLC Registered"
$ws1.Cells.Item(240, 3).Value = 6231501.4
$ws1.Cells.Item(240, 4).Value = 0
$ws1.Cells.Item(240, 5).Value = 0
$ws1.Cells.Item(240, 6).Value = 0
$ws1.Cells.Item(240, 7).Value = 0
$ws1.Cells.Item(240, 8).Value = 6231501.4

$ws1.Cells.Item(241, 1).Value = 46073
$ws1.Cells.Item(241, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(241, 2).Value = "STONEX PRECIOUS METALS LLC Eligible"
$ws1.Cells.Item(241, 3).Value = 1537051.72
$ws1.Cells.Item(241, 4).Value = 0
$ws1.Cells.Item(241, 5).Value = 0
$ws1.Cells.Item(241, 6).Value = 0
$ws1.Cells.Item(241, 7).Value = 0
$ws1.Cells.Item(241, 8).Value = 1537051.72

# ---- Sheet 2: Today_Summary - update Eligible (B) and Total_Stock (D) for depositories with withdrawals ----
$ws2 = $wb.Worksheets.Item("Today_Summary")

$ws2.Cells.Item(4, 2).Value = 13856687.823
$ws2.Cells.Item(4, 4).Value = 26031539.392

$ws2.Cells.Item(5, 2).Value = 15771876.945
$ws2.Cells.Item(5, 4).Value = 17304653.368

$ws2.Cells.Item(6, 2).Value = 19011209.353
$ws2.Cells.Item(6, 4).Value = 22423366.923

$ws2.Cells.Item(8, 2).Value = 142890665.333
$ws2.Cells.Item(8, 4).Value = 154891009.103

# ---- Sheet 3: Monthly_Stats - update monthly cumulative totals ----
$ws3 = $wb.Worksheets.Item("Monthly_Stats")

$ws3.Cells.Item(2, 2).Value = 275812073.613
$ws3.Cells.Item(2, 4).Value = 364003132.877

$ws3.Cells.Item(11, 4).Value = 5145962.733
$ws3.Cells.Item(11, 5).Value = 13856687.823

$ws3.Cells.Item(13, 4).Value = 691419.1240000001
$ws3.Cells.Item(13, 5).Value = 15771876.945

$ws3.Cells.Item(15, 4).Value = 2309846.81
$ws3.Cells.Item(15, 5).Value = 19011209.353

$ws3.Cells.Item(19, 4).Value = 13973136.1
$ws3.Cells.Item(19, 5).Value = 142890665.333
